$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (and correspondingly the <sheet name=.../> entry in workbook.xml)
$ws.Name = "Work Received inc amavat 2021"

# Insert a new first column, shifting the existing A:M data to B:N
$ws.Columns("A:A").Insert()

# Give the new column A (rows 2-6) the same header-row style (bold, bordered,
# centered) already used for the original header cells, without creating a
# duplicate style entry
$ws.Range("B1").Copy()
$ws.Range("A2:A6").PasteSpecial(-4122)

# Populate the new "ranking" column with sequential numbers
$ws.Range("A2").Value = 7
$ws.Range("A3").Value = 8
$ws.Range("A4").Value = 9
$ws.Range("A5").Value = 10
$ws.Range("A6").Value = 11
